$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 67

# Columns A-D hold plain text values (date/time/weekday/week stored as text,
# not native Excel date/time types). Prefix with an apostrophe to force
# text entry (avoids auto-conversion to date/time serials or numbers),
# then reset the style back to Normal so no extra formatting/quote-prefix
# marker is left behind on the cells (matching the rest of the sheet).
$ws.Cells.Item($row, 1).Value = "'2025-02-13"
$ws.Cells.Item($row, 2).Value = "'13:26:25"
$ws.Cells.Item($row, 3).Value = "'Thursday"
$ws.Cells.Item($row, 4).Value = "'06"
$ws.Range("A67:D67").Style = "Normal"

# Columns E-T hold numeric values
$ws.Cells.Item($row, 5).Value = 128315
$ws.Cells.Item($row, 6).Value = 142291
$ws.Cells.Item($row, 7).Value = 169811
$ws.Cells.Item($row, 8).Value = 159055
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 144646
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 191893
$ws.Cells.Item($row, 14).Value = 115099
$ws.Cells.Item($row, 15).Value = 44978
$ws.Cells.Item($row, 16).Value = 28623
$ws.Cells.Item($row, 17).Value = 65377
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 44426
$ws.Cells.Item($row, 20).Value = -1
